# Regenerate save_data to use K (strikeouts) instead of Strike# in column G.
# Only column G ("K") values change for rows 2-10 based on the recalculated
# std/mean and s_vals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 2
    3  = 0
    4  = 1
    5  = 1
    6  = 0
    7  = 1
    9  = 1
    10 = 0
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
